$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Certificate number (appears twice, identical change both times)
Replace-All "2445987" "2409865"

# Date (appears twice, identical change both times)
Replace-All "01.05.2024" "05.04.2024"

# "Светлый" -> "Светлый " (trailing space added, appears twice, identical change both times)
Replace-All "Светлый" "Светлый "

# Representative name/role change
Replace-All "генеральный директор Котлярчук О. Е." "Суперинтендант Мухин К. А."

# Basis document change
Replace-All "Устава" "Доверенности №  от -- "

# Survey / document description change
Replace-All "Рассмотрение технической документации ""Грузовая марка"" № 5234-234234-23 на т/х ""СИНЕГОРСК"" РС 021026" "Ежегодное освидетельствование ССП № 24.42.03.00765.121 от 02.05.2021 "

# Approval letter change
Replace-All "Письмо об одобрении № 121-212-08-343489 от 05.05.2024" "Свидетельство ф. 7.1.27 № 24.02.42.00987.121 от --"

# Monetary amounts
Replace-All "5 150,00 p. (пять тысяч сто пятьдесят рублей 00 копеек)" "15 600,00 p. (пятнадцать тысяч шестьсот рублей 00 копеек)"
Replace-All "1 030,00 p. (одна тысяча тридцать рублей 00 копеек)" "3 120,00 p. (три тысячи сто двадцать рублей 00 копеек)"
Replace-All "6 180,00 p. (шесть тысяч сто восемьдесят рублей 00 копеек)" "18 720,00 p. (восемнадцать тысяч семьсот двадцать рублей 00 копеек)"

# Signature name
Replace-All "О. Е. Котлярчук" "К. А. Мухин"
